$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the "Meta description: ..." paragraph that currently
# sits right under the title heading.
# ---------------------------------------------------------------------
$metaPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Meta description*") {
        $metaPara = $p
        break
    }
}
if ($metaPara -ne $null) {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# Step 2: locate the paragraph "Free spins mode challenging to
# achieve" - the new paragraph goes right after it (i.e. right before
# the final, last paragraph of the document).
# ---------------------------------------------------------------------
$priorPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Free spins mode challenging to achieve*") {
        $priorPara = $p
        break
    }
}

# ---------------------------------------------------------------------
# Step 3: insert a brand-new paragraph right after it, bearing the
# bold text that used to be the meta description's title.
# ---------------------------------------------------------------------
$insertionPoint = $priorPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$count = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($count - 1)

$newPara.Style = "Normal"
$newRange = $newPara.Range
$newRange.Text = "Play Free Elephant King Online Slot Game Review"
$newRange.Font.Reset()

$newTextOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newTextOnly.Font.Bold = 1

# ---------------------------------------------------------------------
# Step 4: replace the text of the final paragraph (the italic prompt)
# with the new meta-description copy, keeping the italic formatting.
# ---------------------------------------------------------------------
$finalPara = $d.Paragraphs.Last
$finalTextOnly = $d.Range($finalPara.Range.Start, $finalPara.Range.End - 1)
$finalTextOnly.Text = "Explore the beauty of African savanna with the Elephant King online slot game. Play for free and experience high potential for payouts with up to 40 winning lines."
